# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# Swap the match-data (columns B through AD) between specific pairs of rows.
# Column A (the sequential row id) stays untouched on each row; only the
# underlying match record (match id, teams, odds, results, etc.) moves
# between the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of row numbers whose B:AD contents must be exchanged.
$pairs = @(
    @(13, 15),
    @(47, 48),
    @(66, 67),
    @(115, 116),
    @(118, 119),
    @(193, 194)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $data1 = $range1.Value2
    $data2 = $range2.Value2

    $range1.Value2 = $data2
    $range2.Value2 = $data1
}
